$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (the "Förändrad" column) for rows 2 through 43:
# value 45777 (date serial) -> 45778
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45777) {
        $cell.Value2 = 45778
    }
}
